$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1966.1666
$ws.Cells.Item(80, 10).Value = 2699.875
$ws.Cells.Item(80, 12).Value = 8099.625
$ws.Cells.Item(80, 14).Value = -10095.625

$ws.Cells.Item(83, 8).Value = 1966.1666
$ws.Cells.Item(83, 10).Value = 2699.875
$ws.Cells.Item(83, 12).Value = 24298.875
$ws.Cells.Item(83, 14).Value = -34282.875

$ws.Cells.Item(132, 8).Value = 1533.2
$ws.Cells.Item(132, 9).Value = 1463.909
$ws.Cells.Item(132, 11).Value = 4391.727000000001
$ws.Cells.Item(132, 13).Value = -1861.727000000001

$ws.Cells.Item(137, 8).Value = 502855.88
$ws.Cells.Item(137, 9).Value = 2245.5
$ws.Cells.Item(137, 11).Value = 6736.5
$ws.Cells.Item(137, 13).Value = -4186.5

$ws.Cells.Item(138, 8).Value = 1411.25
$ws.Cells.Item(138, 9).Value = 1411.25
$ws.Cells.Item(138, 11).Value = 4233.75
$ws.Cells.Item(138, 13).Value = 906.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 36773.83
$ws.Cells.Item(74, 9).Value = 51492.35
$ws.Cells.Item(74, 10).Value = 4066
$ws.Cells.Item(74, 11).Value = 51492.35
$ws.Cells.Item(74, 12).Value = 4066
$ws.Cells.Item(74, 13).Value = -50618.35
$ws.Cells.Item(74, 14).Value = -5814

$ws.Cells.Item(77, 8).Value = 36773.83
$ws.Cells.Item(77, 9).Value = 51492.35
$ws.Cells.Item(77, 10).Value = 4066
$ws.Cells.Item(77, 11).Value = 257461.75
$ws.Cells.Item(77, 12).Value = 20330
$ws.Cells.Item(77, 13).Value = -253093.75
$ws.Cells.Item(77, 14).Value = -29066

$ws.Cells.Item(102, 8).Value = 51808.09
$ws.Cells.Item(102, 9).Value = 60133.35
$ws.Cells.Item(102, 11).Value = 60133.35
$ws.Cells.Item(102, 13).Value = -58511.35

$ws.Cells.Item(124, 8).Value = 62419
$ws.Cells.Item(124, 10).Value = 62419
$ws.Cells.Item(124, 12).Value = 62419
$ws.Cells.Item(124, 14).Value = -72239

$ws.Cells.Item(132, 8).Value = 2540.5
$ws.Cells.Item(132, 9).Value = 2493
$ws.Cells.Item(132, 11).Value = 7479
$ws.Cells.Item(132, 13).Value = -4949

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 169652.83
$ws.Cells.Item(22, 9).Value = 202384.2
$ws.Cells.Item(22, 11).Value = 202384.2
$ws.Cells.Item(22, 13).Value = -202211.2

$ws.Cells.Item(86, 8).Value = 3480.4119
$ws.Cells.Item(86, 9).Value = 3397.7856
$ws.Cells.Item(86, 10).Value = 3866
$ws.Cells.Item(86, 11).Value = 3397.7856
$ws.Cells.Item(86, 12).Value = 3866
$ws.Cells.Item(86, 13).Value = -2274.7856
$ws.Cells.Item(86, 14).Value = -6112

$ws.Cells.Item(89, 8).Value = 3480.4119
$ws.Cells.Item(89, 9).Value = 3397.7856
$ws.Cells.Item(89, 10).Value = 3866
$ws.Cells.Item(89, 11).Value = 16988.928
$ws.Cells.Item(89, 12).Value = 19330
$ws.Cells.Item(89, 13).Value = -11372.928
$ws.Cells.Item(89, 14).Value = -30562

$ws.Cells.Item(94, 8).Value = 2925.1667
$ws.Cells.Item(94, 9).Value = 2010.2
$ws.Cells.Item(94, 10).Value = 7500
$ws.Cells.Item(94, 11).Value = 2010.2
$ws.Cells.Item(94, 12).Value = 7500
$ws.Cells.Item(94, 13).Value = -1559.2
$ws.Cells.Item(94, 14).Value = -8402

$ws.Cells.Item(105, 8).Value = 37188.145
$ws.Cells.Item(105, 9).Value = 44566
$ws.Cells.Item(105, 10).Value = 3250
$ws.Cells.Item(105, 11).Value = 44566
$ws.Cells.Item(105, 12).Value = 3250
$ws.Cells.Item(105, 13).Value = -42819
$ws.Cells.Item(105, 14).Value = -6744

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1538.2727
$ws.Cells.Item(58, 9).Value = 1365.5555
$ws.Cells.Item(58, 11).Value = 1365.5555
$ws.Cells.Item(58, 13).Value = -1162.5555

$ws.Cells.Item(136, 8).Value = 1538.2727
$ws.Cells.Item(136, 9).Value = 1365.5555
$ws.Cells.Item(136, 11).Value = 4096.666499999999
$ws.Cells.Item(136, 13).Value = -1546.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 89.333336
$ws.Cells.Item(17, 9).Value = 19
$ws.Cells.Item(17, 11).Value = 57
$ws.Cells.Item(17, 13).Value = 112

$ws.Cells.Item(34, 8).Value = 2700
$ws.Cells.Item(34, 10).Value = 3303.75
$ws.Cells.Item(34, 12).Value = 9911.25
$ws.Cells.Item(34, 14).Value = -10079.25

$ws.Cells.Item(39, 8).Value = 5936.6
$ws.Cells.Item(39, 9).Value = 400
$ws.Cells.Item(39, 10).Value = 6228
$ws.Cells.Item(39, 11).Value = 1200
$ws.Cells.Item(39, 12).Value = 18684
$ws.Cells.Item(39, 13).Value = -906
$ws.Cells.Item(39, 14).Value = -19272

$ws.Cells.Item(55, 8).Value = 892.2308
$ws.Cells.Item(55, 9).Value = 860
$ws.Cells.Item(55, 10).Value = 999.6667
$ws.Cells.Item(55, 11).Value = 2580
$ws.Cells.Item(55, 12).Value = 2999.0001
$ws.Cells.Item(55, 13).Value = -2403
$ws.Cells.Item(55, 14).Value = -3353.0001

$ws.Cells.Item(68, 8).Value = 126898.875
$ws.Cells.Item(68, 10).Value = 144670.14
$ws.Cells.Item(68, 12).Value = 434010.42
$ws.Cells.Item(68, 14).Value = -435632.42

$ws.Cells.Item(71, 8).Value = 126898.875
$ws.Cells.Item(71, 10).Value = 144670.14
$ws.Cells.Item(71, 12).Value = 1302031.26
$ws.Cells.Item(71, 14).Value = -1310143.26

$ws.Cells.Item(113, 8).Value = 2210828.8
$ws.Cells.Item(113, 10).Value = 2701827
$ws.Cells.Item(113, 12).Value = 8105481
$ws.Cells.Item(113, 14).Value = -8109821

$ws.Cells.Item(129, 8).Value = 37037508
$ws.Cells.Item(129, 9).Value = 531.7143
$ws.Cells.Item(129, 11).Value = 1595.1429
$ws.Cells.Item(129, 13).Value = 3404.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3908.6365
$ws.Cells.Item(126, 9).Value = 2616.6667
$ws.Cells.Item(126, 11).Value = 7850.000100000001
$ws.Cells.Item(126, 13).Value = -5380.000100000001

$ws.Cells.Item(132, 8).Value = 4833.263
$ws.Cells.Item(132, 9).Value = 4152.75
$ws.Cells.Item(132, 10).Value = 5999.857
$ws.Cells.Item(132, 11).Value = 12458.25
$ws.Cells.Item(132, 12).Value = 17999.571
$ws.Cells.Item(132, 13).Value = -9928.25
$ws.Cells.Item(132, 14).Value = -23059.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 9725447
$ws.Cells.Item(40, 9).Value = 2963.182
$ws.Cells.Item(40, 10).Value = 31114912
$ws.Cells.Item(40, 11).Value = 2963.182
$ws.Cells.Item(40, 12).Value = 31114912
$ws.Cells.Item(40, 13).Value = -2827.182
$ws.Cells.Item(40, 14).Value = -31115184

$ws.Cells.Item(93, 8).Value = 1494.5
$ws.Cells.Item(93, 9).Value = 989
$ws.Cells.Item(93, 10).Value = 2000
$ws.Cells.Item(93, 11).Value = 989
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = 259
$ws.Cells.Item(93, 14).Value = -4496

$ws.Cells.Item(134, 8).Value = 137899.33
$ws.Cells.Item(134, 10).Value = 137899.33
$ws.Cells.Item(134, 12).Value = 137899.33
$ws.Cells.Item(134, 14).Value = -148039.33

$ws.Cells.Item(137, 8).Value = 109997.5
$ws.Cells.Item(137, 10).Value = 109997.5
$ws.Cells.Item(137, 12).Value = 109997.5
$ws.Cells.Item(137, 14).Value = -120197.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2452.3704
$ws.Cells.Item(132, 9).Value = 2237
$ws.Cells.Item(132, 11).Value = 6711
$ws.Cells.Item(132, 13).Value = -4181
